$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: Bluetooth / Vcc / 5V, H53 styled like H13 (red fill) ---
$ws.Range("B53").Value = "Bluetooth"
$ws.Range("E53").Value = "Vcc"
$ws.Range("F53").Value = "5V"
$ws.Range("H13").Copy()
$ws.Range("H53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 54: GND / GND, H54 styled like H14 (black/theme fill) ---
$ws.Range("E54").Value = "GND"
$ws.Range("F54").Value = "GND"
$ws.Range("H14").Copy()
$ws.Range("H54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 55/56: TXD/RXD labels first (column E), then pin names (column F) ---
$ws.Range("E55").Value = "TXD"
$ws.Range("E56").Value = "RXD"
$ws.Range("F55").Value = "PE4 (UART05 RX)"
$ws.Range("F56").Value = "PE5 (UART05 TX)"

# --- Row 55: H55 = "TX" with new teal fill + centered ---
$ws.Range("H55").Value = "TX"
$ws.Range("H55").Interior.Color = 15773696
$ws.Range("H55").HorizontalAlignment = -4108
$ws.Range("H55").VerticalAlignment = -4108

# --- Row 56: H56 = "RX" with green fill (same as H47) + centered ---
$ws.Range("H56").Value = "RX"
$ws.Range("H47").Copy()
$ws.Range("H56").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H56").HorizontalAlignment = -4108
$ws.Range("H56").VerticalAlignment = -4108

# --- View bookkeeping: scroll position + active selection ---
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("J52").Select()
